# Aggiornamento 15, 16, 17 marzo
# Append three new daily rows (227-229) to the Polinago report sheet,
# continuing the existing series (dates 44301, 44302, 44303; 0 new
# positive cases; rolling 7-day sum = 1; rolling 7-day sum per 100k = 62.34413965087282),
# reusing the same formatting as the last existing row (226).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style/border/font/number format) of the last
# existing data row onto the three new rows before writing values, so
# that no new cell style gets synthesized and the existing style index
# is reused (matches row 226's date style on column A).
$ws.Range("A226:D226").Copy()
$ws.Range("A227:D229").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$newRows = @(
    @{ Row = 227; Date = 44301; Nuovi = 0; Somma7 = 1; Somma7per100k = 62.34413965087282 },
    @{ Row = 228; Date = 44302; Nuovi = 0; Somma7 = 1; Somma7per100k = 62.34413965087282 },
    @{ Row = 229; Date = 44303; Nuovi = 0; Somma7 = 1; Somma7per100k = 62.34413965087282 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Date
    $ws.Cells.Item($r.Row, 2).Value = $r.Nuovi
    $ws.Cells.Item($r.Row, 3).Value = $r.Somma7
    $ws.Cells.Item($r.Row, 4).Value = $r.Somma7per100k
}

Write-Host "Added rows 227-229 to sheet '$($ws.Name)'"
